$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# The workbook originally has three sheets: "Static", "VF-2 Variables"
# and "ST-10 Variables". "VF-2 Variables" is a superset of
# "ST-10 Variables" (it has 8 extra rows: the "?Q600 1813".."?Q600 1820"
# vibration rows, plus a "?Q600 3006" / "Programmable stop with message"
# row that also exists, unused, on "ST-10 Variables").
#
# The edit merges them into a single "Variable" sheet: the unique
# "VF-2" rows are copied onto "ST-10 Variables" (which keeps its
# sheetId/rId), the duplicate/obsolete "?Q600 3006" row is removed from
# everywhere, "VF-2 Variables" is deleted, and the surviving sheet is
# renamed "Variable" and made the active tab.
# ------------------------------------------------------------------

$wsKeep = $wb.Worksheets.Item("ST-10 Variables")

# Insert the 8 rows that only exist on "VF-2 Variables" (tool 13..20
# vibration readings), right before the existing "?Q600 3004" row.
$wsKeep.Rows.Item(30).Resize(8).Insert()

$extraRows = @(
    @("?Q600 1813", "Max recorded vibrations of tool 13"),
    @("?Q600 1814", "Max recorded vibrations of tool 14"),
    @("?Q600 1815", "Max recorded vibrations of tool 15"),
    @("?Q600 1816", "Max recorded vibrations of tool 16"),
    @("?Q600 1817", "Max recorded vibrations of tool 17"),
    @("?Q600 1818", "Max recorded vibrations of tool 18"),
    @("?Q600 1819", "Max recorded vibrations of tool 19"),
    @("?Q600 1820", "Max recorded vibrations of tool 20")
)
for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $r = 30 + $i
    $wsKeep.Cells.Item($r, 1).Value = $extraRows[$i][0]
    $wsKeep.Cells.Item($r, 2).Value = $extraRows[$i][1]
}

# Drop the obsolete "?Q600 3006" / "Programmable stop with message" row
# (now sitting at row 39 after the insert above).
$wsKeep.Rows.Item(39).Delete()

# Remove the now-redundant "VF-2 Variables" sheet.
$wsDrop = $wb.Worksheets.Item("VF-2 Variables")
$wsDrop.Delete()

# Re-fetch the surviving sheet (the old reference goes stale once the
# workbook's sheet collection changes) and rename it.
$wsKeep = $wb.Worksheets.Item("ST-10 Variables")
$wsKeep.Name = "Variable"

# Update the saved selections/active-tab to match the edited workbook:
# "Static" is no longer the active tab, and its cursor moved to B11.
$wsStatic = $wb.Worksheets.Item("Static")
$wsStatic.Range("B11").Select()

# "Variable" becomes the active tab, scrolled down with row 39 selected.
$wsKeep.Activate()
$wsKeep.Range("A16").Select()
$wsKeep.Rows.Item(39).Select()
